# lab03_Section2.xlsx - "Credits for delayed submissions of lab03"
# Record attendance/credit + notes for a few students who previously had
# no Attendance/Credit entry, and add a late-note for one who already had
# full credit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Battley, Nick (row 5) - attended but didn't demo to the TA
$ws.Range("B5").Value = "Yes"
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = "Attended the lab, but didn't show TA the result."

# Evans, Tyler (row 14) - no credit
$ws.Range("B14").Value = "No"
$ws.Range("C14").Value = 0

# Obannon-Liggins, Joshua (row 23) - no credit
$ws.Range("B23").Value = "No"
$ws.Range("C23").Value = 0

# Siemionek, Oliver (row 30) - no credit
$ws.Range("B30").Value = "No"
$ws.Range("C30").Value = 0

# Wu, Michael (row 35) - already had partial credit, add a note about the delay
$ws.Range("D35").Value = "Delayed submission on Monday."

# Restore the view roughly where the author left it
$ws.Range("H34").Select()
$excel.ActiveWindow.ScrollRow = 20
